$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.362.77"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.844.72"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'240.05"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.6342"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07544"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.2958"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'24.68"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'4.983"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'0.6829"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'83.11"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'0.000009884"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "'6.162"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").Value = "29.406.43"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'230.37"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D21").Value = "'7.541"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +236.65%  "
$ws.Range("E24").Value = "  +170.76%  "
$ws.Range("D25").Value = "'156.24"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'0.1406"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "'17.67"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "'1.469"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'0.05706"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "'1.251"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'4.028"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "'1.853"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("D35").Value = "'1.156"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").Value = "'0.7152"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "'2.597"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "1.249.64"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").Value = "'2.799"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'0.01811"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("E41").Value = "  +267.79%  "
$ws.Range("D42").Value = "'0.9021"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'101.84"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'66.33"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "'7.068"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("D47").Value = "'9.126"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'0.4016"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'1.701"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "'0.05740"
$ws.Range("E51").Value = "  -0.24%  "
